{"js": "// The document contains a single table of two-digit division problems,\n// written as text like \"95\u00f74=23, 3\" (dividend\u00f7divisor=quotient, remainder).\n// Every non-empty cell's text needs to be replaced with a new division\n// problem, in table (row-major) order. Blank rows/cells are left untouched.\nconst replacements = [\n  \"68\u00f79=7, 5\",\n  \"89\u00f79=9, 8\",\n  \"78\u00f76=13, 0\",\n  \"91\u00f76=15, 1\",\n  \"14\u00f79=1, 5\",\n  \"39\u00f72=19, 1\",\n  \"96\u00f78=12, 0\",\n  \"25\u00f73=8, 1\",\n  \"22\u00f74=5, 2\",\n  \"79\u00f74=19, 3\",\n  \"29\u00f77=4, 1\",\n  \"70\u00f76=11, 4\",\n  \"21\u00f74=5, 1\",\n  \"84\u00f73=28, 0\",\n  \"48\u00f72=24, 0\",\n  \"49\u00f77=7, 0\",\n  \"65\u00f74=16, 1\",\n  \"78\u00f78=9, 6\",\n  \"55\u00f73=18, 1\",\n  \"74\u00f73=24, 2\",\n  \"96\u00f73=32, 0\",\n  \"29\u00f72=14, 1\",\n  \"10\u00f74=2, 2\",\n  \"26\u00f73=8, 2\",\n  \"98\u00f79=10, 8\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet repIdx = 0;\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    for (const cell of cells.items) {\n      cell.body.load(\"text\");\n      await context.sync();\n\n      const text = cell.body.text.replace(/[\\r\\x07]+$/, \"\");\n      if (text.trim().length > 0) {\n        if (repIdx >= replacements.length) {\n          throw new Error(\"Ran out of replacement values for non-empty cell: \" + text);\n        }\n        cell.value = replacements[repIdx];\n        repIdx++;\n      }\n    }\n  }\n}\n\nawait context.sync();\n\nif (repIdx !== replacements.length) {\n  throw new Error(`Expected to replace ${replacements.length} cells, replaced ${repIdx}`);\n}\n", "ps1": "# The document contains a single table of two-digit division problems,\n# written as text like \"95\u00f74=23, 3\" (dividend\u00f7divisor=quotient, remainder).\n# Every non-empty cell's text needs to be replaced with a new division\n# problem, in table (row-major) order. Blank rows/cells are left untouched.\n$replacements = @(\n  \"68\u00f79=7, 5\",\n  \"89\u00f79=9, 8\",\n  \"78\u00f76=13, 0\",\n  \"91\u00f76=15, 1\",\n  \"14\u00f79=1, 5\",\n  \"39\u00f72=19, 1\",\n  \"96\u00f78=12, 0\",\n  \"25\u00f73=8, 1\",\n  \"22\u00f74=5, 2\",\n  \"79\u00f74=19, 3\",\n  \"29\u00f77=4, 1\",\n  \"70\u00f76=11, 4\",\n  \"21\u00f74=5, 1\",\n  \"84\u00f73=28, 0\",\n  \"48\u00f72=24, 0\",\n  \"49\u00f77=7, 0\",\n  \"65\u00f74=16, 1\",\n  \"78\u00f78=9, 6\",\n  \"55\u00f73=18, 1\",\n  \"74\u00f73=24, 2\",\n  \"96\u00f73=32, 0\",\n  \"29\u00f72=14, 1\",\n  \"10\u00f74=2, 2\",\n  \"26\u00f73=8, 2\",\n  \"98\u00f79=10, 8\"\n)\n\n$d = $word.ActiveDocument\n$repIdx = 0\n\nforeach ($t in $d.Tables) {\n  for ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n      $cell = $t.Cell($r, $c)\n      $raw = $cell.Range.Text\n      $trimmed = $raw.TrimEnd([char]13, [char]7).Trim()\n      if ($trimmed.Length -gt 0) {\n        $cell.Range.Text = $replacements[$repIdx]\n        $repIdx = $repIdx + 1\n      }\n    }\n  }\n}\n\nWrite-Output (\"Replaced \" + $repIdx + \" cells\")\n"}
